# Insert a new weekly record as row 142 in the "Uva" (grape) price sheet.
# This pushes the existing rows 142..168 down to 143..169 (the row that was
# previously last, row 168, becomes row 169), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 142, shifting rows below down.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record's data.
$ws.Cells.Item(142, 1).Value  = 8
$ws.Cells.Item(142, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(142, 3).Value  = "Coquimbo"
$ws.Cells.Item(142, 4).Value  = 45275
$ws.Cells.Item(142, 5).Value  = 4
$ws.Cells.Item(142, 6).Value  = "Fruta"
$ws.Cells.Item(142, 7).Value  = 100109
$ws.Cells.Item(142, 8).Value  = "Uva"
$ws.Cells.Item(142, 9).Value  = 100109001
$ws.Cells.Item(142, 10).Value = "Uva"
$ws.Cells.Item(142, 11).Value = "Flame Seedless"
$ws.Cells.Item(142, 12).Value = "Primera"
$ws.Cells.Item(142, 13).Value = 700
$ws.Cells.Item(142, 14).Value = 10000
$ws.Cells.Item(142, 15).Value = 11000
$ws.Cells.Item(142, 16).Value = 10500
$ws.Cells.Item(142, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(142, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(142, 19).Value = 1050
$ws.Cells.Item(142, 20).Value = 10
